$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the User Name cell (A2) from "Juliana " to "Regiane"
$ws.Range("A2").Value = "Regiane"

# Move the active selection from D2 to A2
$ws.Range("A2").Select()
